$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated run timestamps in column C (test re-run a bit later than before)
$ws.Range("C2").Value = '01/12/2025 19:31:05'
$ws.Range("C3").Value = '01/12/2025 19:31:09'
$ws.Range("C4").Value = '01/12/2025 19:31:11'
$ws.Range("C5").Value = '01/12/2025 19:31:16'
$ws.Range("C6").Value = '01/12/2025 19:31:21'
$ws.Range("C7").Value = '01/12/2025 19:31:21'
$ws.Range("C8").Value = '01/12/2025 19:31:40'
$ws.Range("C9").Value = '01/12/2025 19:31:48'
$ws.Range("C10").Value = '01/12/2025 19:31:53'
$ws.Range("C11").Value = '01/12/2025 19:32:02'
$ws.Range("C12").Value = '01/12/2025 19:32:06'
$ws.Range("C13").Value = '01/12/2025 19:32:08'
$ws.Range("C14").Value = '01/12/2025 19:32:13'
$ws.Range("C15").Value = '01/12/2025 19:32:18'
$ws.Range("C16").Value = '01/12/2025 19:32:18'
$ws.Range("C17").Value = '01/12/2025 19:32:37'
$ws.Range("C18").Value = '01/12/2025 19:32:40'
$ws.Range("C19").Value = '01/12/2025 19:32:45'
$ws.Range("D19").Value = '✓ המערכת זיהתה אימייל קיים והציגה שגיאה: "יש לבחור יישוב ורחוב." - הבדיקה עברה בהצלחה!'
$ws.Range("C20").Value = '01/12/2025 19:33:00'
$ws.Range("C21").Value = '01/12/2025 19:33:04'
$ws.Range("C22").Value = '01/12/2025 19:33:06'
$ws.Range("C23").Value = '01/12/2025 19:33:11'
$ws.Range("C24").Value = '01/12/2025 19:33:16'
$ws.Range("C25").Value = '01/12/2025 19:33:16'
$ws.Range("C26").Value = '01/12/2025 19:33:34'
$ws.Range("C27").Value = '01/12/2025 19:33:38'
$ws.Range("C28").Value = '01/12/2025 19:33:42'
$ws.Range("C29").Value = '01/12/2025 19:34:00'
$ws.Range("C30").Value = '01/12/2025 19:34:04'
$ws.Range("C31").Value = '01/12/2025 19:34:06'
$ws.Range("C32").Value = '01/12/2025 19:34:11'
$ws.Range("C33").Value = '01/12/2025 19:34:16'
$ws.Range("C34").Value = '01/12/2025 19:34:16'
$ws.Range("C35").Value = '01/12/2025 19:34:35'
$ws.Range("C36").Value = '01/12/2025 19:34:38'
$ws.Range("C37").Value = '01/12/2025 19:34:43'

# Column D width shrank because the long "email already exists" error message
# (row 19) was replaced with a shorter validation message; re-fit the column.
$ws.Columns.Item(4).ColumnWidth = 110.28
